$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(7, 12, 54, 55, 75, 80, 99, 103, 121, 126, 145, 172, 192, 196, 213, 218, 237, 241, 259, 264, 285, 288, 337, 338, 357, 361)

foreach ($r in $rows) {
    $cell = $ws.Range("E$r")
    $current = $cell.Value2
    if ($current -eq "Nada") {
        $cell.Value = "Trabaja"
    } elseif ($current -eq "Trabaja") {
        $cell.Value = "Nada"
    }
}
